# Fix onboarding state logic
# The Tasks sheet's columns D/E/F were swapped (priority/createdAt/dueDate
# instead of dueDate/priority/createdAt) and a newly-created onboarding task
# ("jnvjdnvds") was missing from the list. Reorder the columns and insert the
# missing task as the new row 2 (it has no dueDate yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new task at row 2; existing rows 2-3 shift down to 3-4.
$ws.Rows.Item(2).Insert()

# --- Header row: D/E/F now read priority, createdAt, dueDate ---
$ws.Range("D1").Value = "priority"
$ws.Range("E1").Value = "createdAt"
$ws.Range("F1").Value = "dueDate"

# --- Row 2: brand new onboarding task (no dueDate yet) ---
$ws.Range("A2").Value = "00037519-1bcb-4bb8-b2e0-1dff1cb4eb99"
$ws.Range("B2").Value = "jnvjdnvds"
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = "medium"
$ws.Range("E2").Value = 1763735225324

# --- Row 3 (was row 2): re-map its D/E/F values into the new column order ---
$ws.Range("D3").Value = "medium"
$ws.Range("E3").Value = 1763734626268
$f3 = $ws.Range("F3")
$f3.NumberFormat = "@"
$f3.Value = "2025-11-22"
$f3.Style = "Normal"

# --- Row 4 (was row 3): re-map its D/E/F values into the new column order ---
$ws.Range("D4").Value = "medium"
$ws.Range("E4").Value = 1763733656091
$f4 = $ws.Range("F4")
$f4.NumberFormat = "@"
$f4.Value = "2025-11-15"
$f4.Style = "Normal"
